$wb = $excel.ActiveWorkbook

$wsRoute = $wb.Worksheets.Item("Route")
$wsMobil = $wb.Worksheets.Item("Mobil")

# --- Route sheet ---
# Header date cell L3 was a date serial; now it's a literal text date "29/3/2025"
$wsRoute.Range("L3").Value = "29/3/2025"

# Update the view: selection moves to L4:M4 (this also resets the scrolled
# top-left cell back to A1, matching the saved view state)
$wsRoute.Activate()
$wsRoute.Range("L4:M4").Select()

# --- Mobil sheet ---
# Report date at top right
$wsMobil.Range("G1").Value = "2/1/2025"

# Row 5 name changed
$wsMobil.Range("C5").Value = "Chiranjit Barai"

# Updated monthly amounts
$wsMobil.Range("D7").Value = 550
$wsMobil.Range("D9").Value = 550
$wsMobil.Range("D10").Value = 550

$wb.Save()
